# Auto-generated Excel COM-interop script to apply scheduled runner value updates
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 375.16666
$ws.Range("I5").Value = 23.333334
$ws.Range("J5").Value = 727
$ws.Range("K5").Value = 23.333334
$ws.Range("L5").Value = 727
$ws.Range("M5").Value = 91.66666599999999
$ws.Range("N5").Value = -957
$ws.Range("H17").Value = 1200.3684
$ws.Range("J17").Value = 1200.3684
$ws.Range("L17").Value = 3601.1052
$ws.Range("N17").Value = -3937.1052
$ws.Range("H18").Value = 1026.0526
$ws.Range("I18").Value = 776.17645
$ws.Range("J18").Value = 3150
$ws.Range("K18").Value = 776.17645
$ws.Range("L18").Value = 3150
$ws.Range("M18").Value = -492.17645
$ws.Range("N18").Value = -3718
$ws.Range("H19").Value = 947.8570999999999
$ws.Range("I19").Value = 196.25
$ws.Range("J19").Value = 1950
$ws.Range("K19").Value = 196.25
$ws.Range("L19").Value = 1950
$ws.Range("M19").Value = -21.25
$ws.Range("N19").Value = -2300
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
$ws.Range("H40").Value = 60319.234
$ws.Range("I40").Value = 112532.664
$ws.Range("K40").Value = 112532.664
$ws.Range("M40").Value = -112357.664
$ws.Range("H41").Value = 793.6111
$ws.Range("I41").Value = 734.4286
$ws.Range("J41").Value = 831.2727
$ws.Range("K41").Value = 734.4286
$ws.Range("L41").Value = 831.2727
$ws.Range("M41").Value = -294.4286
$ws.Range("N41").Value = -1711.2727
$ws.Range("H43").Value = 2980.1667
$ws.Range("I43").Value = 4660.3335
$ws.Range("J43").Value = 1300
$ws.Range("K43").Value = 4660.3335
$ws.Range("L43").Value = 1300
$ws.Range("M43").Value = -4591.3335
$ws.Range("N43").Value = -1438
$ws.Range("H51").Value = 10378.5
$ws.Range("J51").Value = 3539.9
$ws.Range("L51").Value = 3539.9
$ws.Range("N51").Value = -4507.9
$ws.Range("H55").Value = 233.16667
$ws.Range("J55").Value = 282
$ws.Range("L55").Value = 282
$ws.Range("N55").Value = -710
$ws.Range("H113").Value = 2299.8333
$ws.Range("I113").Value = 2575
$ws.Range("J113").Value = 1749.5
$ws.Range("K113").Value = 2575
$ws.Range("L113").Value = 1749.5
$ws.Range("M113").Value = 679
$ws.Range("N113").Value = -8257.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 200822.2
$ws.Range("I2").Value = 1055.5
$ws.Range("K2").Value = 1055.5
$ws.Range("M2").Value = -942.5
$ws.Range("H5").Value = 60
$ws.Range("I5").Value = 20
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 20
$ws.Range("L5").Value = 100
$ws.Range("M5").Value = 92
$ws.Range("N5").Value = -324
$ws.Range("H45").Value = 92743.55
$ws.Range("I45").Value = 144314.28
$ws.Range("J45").Value = 2494.75
$ws.Range("K45").Value = 144314.28
$ws.Range("L45").Value = 2494.75
$ws.Range("M45").Value = -143937.28
$ws.Range("N45").Value = -3248.75
$ws.Range("H102").Value = 127167.5
$ws.Range("I102").Value = 335726.66
$ws.Range("J102").Value = 2032
$ws.Range("K102").Value = 335726.66
$ws.Range("L102").Value = 2032
$ws.Range("M102").Value = -334104.66
$ws.Range("N102").Value = -5276
$ws.Range("H116").Value = 200822.2
$ws.Range("I116").Value = 1055.5
$ws.Range("K116").Value = 1055.5
$ws.Range("M116").Value = 1238.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 200822.2
$ws.Range("I3").Value = 1055.5
$ws.Range("K3").Value = 1055.5
$ws.Range("M3").Value = -941.5
$ws.Range("H4").Value = 60
$ws.Range("I4").Value = 20
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 20
$ws.Range("L4").Value = 100
$ws.Range("M4").Value = 95
$ws.Range("N4").Value = -330
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H107").Value = 250115820
$ws.Range("I107").Value = 333486100
$ws.Range("K107").Value = 333486100
$ws.Range("M107").Value = -333484180

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 866.6667
$ws.Range("J8").Value = 900
$ws.Range("L8").Value = 900
$ws.Range("N8").Value = -1180
$ws.Range("H15").Value = 15000
$ws.Range("J15").Value = 15000
$ws.Range("L15").Value = 15000
$ws.Range("N15").Value = -15340
$ws.Range("H25").Value = 16600
$ws.Range("I25").Value = 13000
$ws.Range("J25").Value = 19000
$ws.Range("K25").Value = 13000
$ws.Range("L25").Value = 19000
$ws.Range("M25").Value = -12826
$ws.Range("N25").Value = -19348
$ws.Range("H29").Value = 6999.6665
$ws.Range("J29").Value = 6999.6665
$ws.Range("L29").Value = 6999.6665
$ws.Range("N29").Value = -7585.6665
$ws.Range("H47").Value = 35000
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 35000
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 35000
$ws.Range("M47").ClearContents()
$ws.Range("N47").Value = -36132
$ws.Range("H122").Value = 1000
$ws.Range("I122").Value = 1000
$ws.Range("K122").Value = 3000
$ws.Range("M122").Value = -550
$ws.Range("H132").Value = 3712.7407
$ws.Range("I132").Value = 2788.4546
$ws.Range("J132").Value = 7779.6
$ws.Range("K132").Value = 8365.363799999999
$ws.Range("L132").Value = 23338.8
$ws.Range("M132").Value = -5835.363799999999
$ws.Range("N132").Value = -28398.8

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 228.875
$ws.Range("I15").Value = 100
$ws.Range("J15").Value = 271.83334
$ws.Range("K15").Value = 300
$ws.Range("L15").Value = 815.5000200000001
$ws.Range("M15").Value = -160
$ws.Range("N15").Value = -1095.50002
$ws.Range("H117").Value = 3312.8
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 3312.8
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 9938.400000000001
$ws.Range("M117").ClearContents()
$ws.Range("N117").Value = -16822.4
$ws.Range("H131").Value = 752.91
$ws.Range("J131").Value = 780.117
$ws.Range("L131").Value = 2340.351
$ws.Range("N131").Value = -12420.351

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 403654.72
$ws.Range("I102").Value = 3026.1
$ws.Range("K102").Value = 3026.1
$ws.Range("M102").Value = -1404.1

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3099.45
$ws.Range("I7").Value = 2865.72
$ws.Range("K7").Value = 2865.72
$ws.Range("M7").Value = -2753.72
$ws.Range("H9").Value = 300
$ws.Range("I9").Value = 510
$ws.Range("J9").Value = 90
$ws.Range("K9").Value = 510
$ws.Range("L9").Value = 90
$ws.Range("M9").Value = -286
$ws.Range("N9").Value = -538
$ws.Range("H16").Value = 4858364.5
$ws.Range("I16").Value = 6300748
$ws.Range("K16").Value = 6300748
$ws.Range("M16").Value = -6300578
$ws.Range("H22").Value = 2036
$ws.Range("I22").Value = 1539.8572
$ws.Range("K22").Value = 1539.8572
$ws.Range("M22").Value = -1244.8572
$ws.Range("H27").Value = 2036
$ws.Range("I27").Value = 1539.8572
$ws.Range("K27").Value = 1539.8572
$ws.Range("M27").Value = -1432.8572
$ws.Range("H40").Value = 79053.84
$ws.Range("I40").Value = 1000000
$ws.Range("J40").Value = 2308.3333
$ws.Range("K40").Value = 1000000
$ws.Range("L40").Value = 2308.3333
$ws.Range("M40").Value = -999864
$ws.Range("N40").Value = -2580.3333
$ws.Range("H46").Value = 675287.9
$ws.Range("I46").Value = 484.75
$ws.Range("J46").Value = 920670.8
$ws.Range("K46").Value = 484.75
$ws.Range("L46").Value = 920670.8
$ws.Range("M46").Value = -296.75
$ws.Range("N46").Value = -921046.8
$ws.Range("H55").Value = 455596.88
$ws.Range("J55").Value = 805.5
$ws.Range("L55").Value = 805.5
$ws.Range("N55").Value = -1151.5
$ws.Range("H93").Value = 1263.2325
$ws.Range("I93").Value = 1177
$ws.Range("K93").Value = 1177
$ws.Range("M93").Value = 71
$ws.Range("H126").Value = 3099.45
$ws.Range("I126").Value = 2865.72
$ws.Range("K126").Value = 8597.16
$ws.Range("M126").Value = -6127.16

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2667.0625
$ws.Range("I122").Value = 1996.8
$ws.Range("J122").Value = 2971.7273
$ws.Range("K122").Value = 5990.4
$ws.Range("L122").Value = 8915.1819
$ws.Range("M122").Value = -3540.4
$ws.Range("N122").Value = -13815.1819
$ws.Range("H126").Value = 2024.9166
$ws.Range("I126").Value = 1953.2222
$ws.Range("J126").Value = 2240
$ws.Range("K126").Value = 5859.6666
$ws.Range("L126").Value = 6720
$ws.Range("M126").Value = -3389.6666
$ws.Range("N126").Value = -11660
